# Update "想去人数" (want-to-go count) values in column F across the four
# worksheets of the workbook, as published with the gh-pages output at 456a3b4.

$wb = $excel.ActiveWorkbook

# Sheet 1: 展览
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("F3").Value  = 4775
$ws1.Range("F4").Value  = 626
$ws1.Range("F5").Value  = 217
$ws1.Range("F6").Value  = 1925
$ws1.Range("F8").Value  = 816
$ws1.Range("F9").Value  = 41
$ws1.Range("F10").Value = 22
$ws1.Range("F12").Value = 1168
$ws1.Range("F13").Value = 1622
$ws1.Range("F14").Value = 844
$ws1.Range("F15").Value = 517
$ws1.Range("F16").Value = 1979
$ws1.Range("F17").Value = 618
$ws1.Range("F21").Value = 238
$ws1.Range("F22").Value = 102
$ws1.Range("F23").Value = 102
$ws1.Range("F24").Value = 1566
$ws1.Range("F26").Value = 637
$ws1.Range("F27").Value = 2546
$ws1.Range("F31").Value = 1645
$ws1.Range("F35").Value = 75
$ws1.Range("F36").Value = 4380

# Sheet 2: 演出
$ws2 = $wb.Worksheets.Item(2)
$ws2.Range("F7").Value  = 1
$ws2.Range("F8").Value  = 4174
$ws2.Range("F11").Value = 47
$ws2.Range("F12").Value = 3
$ws2.Range("F17").Value = 296
$ws2.Range("F24").Value = 55
$ws2.Range("F29").Value = 188
$ws2.Range("F32").Value = 27

# Sheet 3: 本地生活
$ws3 = $wb.Worksheets.Item(3)
$ws3.Range("F5").Value = 1756
$ws3.Range("F7").Value = 415

# Sheet 4: 全部类型
$ws4 = $wb.Worksheets.Item(4)
$ws4.Range("F4").Value  = 1756
$ws4.Range("F6").Value  = 415
$ws4.Range("F9").Value  = 4775
$ws4.Range("F10").Value = 217
$ws4.Range("F11").Value = 1925
$ws4.Range("F15").Value = 22
$ws4.Range("F17").Value = 1168
$ws4.Range("F18").Value = 1622
$ws4.Range("F20").Value = 47
$ws4.Range("F22").Value = 844
$ws4.Range("F23").Value = 517
$ws4.Range("F24").Value = 1979
$ws4.Range("F25").Value = 618
$ws4.Range("F29").Value = 238
$ws4.Range("F31").Value = 102
$ws4.Range("F32").Value = 102
$ws4.Range("F34").Value = 296
$ws4.Range("F36").Value = 1566
$ws4.Range("F38").Value = 637
$ws4.Range("F41").Value = 2546
$ws4.Range("F45").Value = 1645
$ws4.Range("F50").Value = 4380
